# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
# Mirrors the existing header formatting (bold, thin box border,
# centered horizontally / top-aligned vertically) used by row 1's
# other header cells, and fills every data row (2-46) with the
# team's 2013 season record: 96 wins, 66 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- Data rows (2 through 46) ------------------------------------------
$lastRow = 46
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 96   # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 66   # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF -> Ties
}
